$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume values scraped on Tue Apr 18 13:33:27 UTC 2023
# (includes row 14/15 and row 42/43 coin re-ranking swaps)

$ws.Range("D2").Value = "30.527.97"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").Value = "2.124.89"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'347.22"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.5244"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").Value = "'0.4476"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "'54.40"
$ws.Range("E9").Value = "  +5.15%  "
$ws.Range("D10").Value = "'0.09416"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "'8.702"
$ws.Range("E13").Value = "  +7.05%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.126.06"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'6.968"
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "'102.37"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "'0.06735"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'6.353"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").Value = "'1.006"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "30.541.24"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "2.359.32"
$ws.Range("E26").Value = "  +1.14%  "
$ws.Range("D27").Value = "'22.27"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "'2.554"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "'162.46"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'134.53"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "'1.165"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").Value = "'1.779"
$ws.Range("E32").Value = "  +9.40%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "'6.879"
$ws.Range("E34").Value = "  +12.32%  "
$ws.Range("D35").Value = "'6.301"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("D36").Value = "'3.965"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "'10.65"
$ws.Range("E37").Value = "  +5.13%  "
$ws.Range("D38").Value = "'0.02662"
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("D39").Value = "'0.06874"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "'0.7145"
$ws.Range("E40").Value = "  +4.24%  "
$ws.Range("D41").Value = "'12.69"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.341"
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.2240"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").Value = "'0.6954"
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("D45").Value = "'14.68"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").Value = "'2.393"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'1.327"
$ws.Range("E48").Value = "  +14.10%  "
$ws.Range("D49").Value = "'3.655"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "'0.00000000346"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +1.28%  "
